# Insert a new weekly data row at row 52, pushing existing rows 52..130 down to 53..131.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new weekly record.
$ws.Cells.Item(52, 1).Value = 9
$ws.Cells.Item(52, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(52, 3).Value = "Metropolitana"
$ws.Cells.Item(52, 4).Value = 44799
$ws.Cells.Item(52, 5).Value = 13
$ws.Cells.Item(52, 6).Value = 100112022
$ws.Cells.Item(52, 7).Value = "Arveja Verde"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 20
$ws.Cells.Item(52, 11).Value = 41000
$ws.Cells.Item(52, 12).Value = 41000
$ws.Cells.Item(52, 13).Value = 41000
$ws.Cells.Item(52, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(52, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(52, 16).Value = 1640
$ws.Cells.Item(52, 17).Value = 25
$ws.Cells.Item(52, 18).Value = "Hortaliza"
